$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency price and volume data

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "35.990.82"
$ws.Range("E2").Value = "  -1.18%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.994.93"
$ws.Range("E3").Value = "  -2.27%  "

$ws.Range("E4").Value = "  +0.00%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "245.67"
$ws.Range("E5").Value = "  +0.20%  "

$ws.Range("E6").Value = "  -2.67%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "59.59"
$ws.Range("E7").Value = "  +11.62%  "

$ws.Range("E8").Value = "  +0.00%  "

$ws.Range("E9").Value = "  -5.68%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.367"
$ws.Range("E10").Value = "  +1.87%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0743"
$ws.Range("E11").Value = "  +0.32%  "

$ws.Range("E12").Value = "  -1.60%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.945"
$ws.Range("E13").Value = "  +2.37%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "14.82"
$ws.Range("E14").Value = "  +2.82%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.281.96"
$ws.Range("E15").Value = "  -2.61%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.41"
$ws.Range("E16").Value = "  +1.12%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "19.30"
$ws.Range("E17").Value = "  +14.34%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.980.94"
$ws.Range("E18").Value = "  -3.00%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "35.900.92"
$ws.Range("E19").Value = "  -1.19%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "71.88"
$ws.Range("E20").Value = "  +1.15%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0₃0852"
$ws.Range("E21").Value = "  +0.43%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.24"
$ws.Range("E22").Value = "  +2.16%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "233.44"
$ws.Range("E23").Value = "  -0.90%  "

$ws.Range("E24").Value = "  +0.07%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.62"
$ws.Range("E25").Value = "  +17.73%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.28"
$ws.Range("E26").Value = "  -3.54%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.68"
$ws.Range("E27").Value = "  +7.14%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "165.86"
$ws.Range("E28").Value = "  +1.44%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "19.55"
$ws.Range("E29").Value = "  -0.91%  "

$ws.Range("E30").Value = "  -0.30%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "5.08"
$ws.Range("E31").Value = "  +2.26%  "

$ws.Range("E32").Value = "  -1.71%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0994"
$ws.Range("E33").Value = "  +15.12%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0604"
$ws.Range("E34").Value = "  +3.16%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.47"
$ws.Range("E35").Value = "  +2.86%  "

$ws.Range("E36").Value = "  +11.93%  "

$ws.Range("E37").Value = "  -0.05%  "

$ws.Range("E38").Value = "  -1.77%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.72"
$ws.Range("E39").Value = "  +15.90%  "

$ws.Range("E40").Value = "  +1.31%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0966"
$ws.Range("E41").Value = "  +8.77%  "

$ws.Range("E42").Value = "  -0.41%  "

$ws.Range("E43").Value = "  +1.55%  "

$ws.Range("B44").Value = "InjectiveProtocol"
$ws.Range("C44").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "16.80"
$ws.Range("E44").Value = "  +8.25%  "

$ws.Range("B45").Value = "ARBITRUM"
$ws.Range("C45").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.11"
$ws.Range("E45").Value = "  +2.36%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "93.87"
$ws.Range("E46").Value = "  +1.33%  "

$ws.Range("E47").Value = "  +5.80%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.367.54"
$ws.Range("E48").Value = "  -0.69%  "

$ws.Range("E49").Value = "  -0.32%  "

$ws.Range("E50").Value = "  +4.02%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "46.52"
$ws.Range("E51").Value = "  +3.92%  "
